# 2016_19_2.xlsx — "Mise à jour du TODO. Correction de bug dans la validation
# du panier. Séparation des états de factures des producteurs"
#
# Semantic changes applied:
#  1. Address string on the letterhead collapsed to a single line
#     ("Chemin de Saint Clair\n07000 PRIVAS" -> "Chemin de Saint Clair, 07000 PRIVAS").
#  2. "Edité le :" timestamp updated to the new export run.
#  3. The "Pain complet" product line (row 16) is removed from the basket
#     entirely — the bug fix to basket validation — which shifts every
#     subsequent row up by one (TOTAL row + footer row included).
#  4. The quantities or the five remaining basket lines are corrected, and
#     the basket is re-sorted (Radis now precedes Salade).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Collapse the two-line address into one line.
$ws.Range("A4").Value = "Chemin de Saint Clair, 07000 PRIVAS"

# 2) Refresh the "édité le" export timestamp.
$ws.Range("B11").Value = "12/05/2016 23:32:20"

# 3) Drop the "Pain complet" row — everything below slides up one row,
#    formulas (D*E, SUBTOTAL) and the dimension/footer follow automatically.
$ws.Rows("16").Delete()

# 4) Fix up the remaining basket rows (now at 16-19) with the corrected
#    quantities, and re-order Radis/Salade.
#    Row 16: Tomates grappe / Légumes / Au poids
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 1

#    Row 17: Pomme de terre / Légumes / Au poids
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 1

#    Row 18: Radis / Légumes / A la pièce
$ws.Range("A18").Value = "Radis"
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 1

#    Row 19 (last basket row, keeps the heavier bottom border): Salade / Légumes / A la pièce
$ws.Range("A19").Value = "Salade"
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 1
